$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: التكليف (assignment) grades for lesson 10
$ws.Range("C23").Value = "التكليف: 5"
$ws.Range("I23").Value = "التكليف: 5"
$ws.Range("L23").Value = "التكليف: 5"
$ws.Range("M23").Value = "التكليف: 5"

# Row 24: الحضور (attendance) grades for lesson 11, with lesson title in N24
$ws.Range("C24").Value = "الحضور: 5"
$ws.Range("D24").Value = "الحضور:5"
$ws.Range("I24").Value = "الحضور: 5"
$ws.Range("K24").Value = "الحضور: 5"
$ws.Range("L24").Value = "الحضور: 5"
$ws.Range("M24").Value = "الحضور: 5"
$ws.Range("N24").Value = "الدرس 11 (آداب فتية حول الرسول)"

# K23 set after N24 so the new shared-string ordering matches (lesson title before "التكليف:5")
$ws.Range("K23").Value = "التكليف:5"

# Row 25: التكليف (assignment) grades for lesson 11 (partial)
$ws.Range("K25").Value = "التكليف: 5"
$ws.Range("M25").Value = "التكليف: 5"

# Update the active selection to reflect where the user left off after entering the new grades
$ws.Range("L26").Select()
